$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet - column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 118
$ws1.Cells.Item(3, 6).Value = 307
$ws1.Cells.Item(4, 6).Value = 185
$ws1.Cells.Item(5, 6).Value = 1178
$ws1.Cells.Item(6, 6).Value = 417
$ws1.Cells.Item(7, 6).Value = 96
$ws1.Cells.Item(8, 6).Value = 126
$ws1.Cells.Item(9, 6).Value = 132
$ws1.Cells.Item(11, 6).Value = 255
$ws1.Cells.Item(12, 6).Value = 150
$ws1.Cells.Item(13, 6).Value = 151
$ws1.Cells.Item(14, 6).Value = 1363
$ws1.Cells.Item(15, 6).Value = 517
$ws1.Cells.Item(16, 6).Value = 196
$ws1.Cells.Item(17, 6).Value = 308
$ws1.Cells.Item(19, 6).Value = 711
$ws1.Cells.Item(20, 6).Value = 1098
$ws1.Cells.Item(21, 6).Value = 56
$ws1.Cells.Item(23, 6).Value = 2547
$ws1.Cells.Item(24, 6).Value = 1306
$ws1.Cells.Item(26, 6).Value = 219
$ws1.Cells.Item(27, 6).Value = 375
$ws1.Cells.Item(28, 6).Value = 944
$ws1.Cells.Item(29, 6).Value = 779
$ws1.Cells.Item(30, 6).Value = 1091
$ws1.Cells.Item(31, 6).Value = 132
$ws1.Cells.Item(32, 6).Value = 88
$ws1.Cells.Item(34, 6).Value = 437
$ws1.Cells.Item(35, 6).Value = 596
$ws1.Cells.Item(36, 6).Value = 762
$ws1.Cells.Item(37, 6).Value = 326
$ws1.Cells.Item(38, 6).Value = 216

# 演出 (Show) sheet - column F ("想去人数")
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(13, 6).Value = 526
$ws2.Cells.Item(21, 6).Value = 12

# 全部类型 (All types) sheet - column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5, 6).Value = 118
$ws4.Cells.Item(6, 6).Value = 307
$ws4.Cells.Item(7, 6).Value = 185
$ws4.Cells.Item(10, 6).Value = 1178
$ws4.Cells.Item(11, 6).Value = 417
$ws4.Cells.Item(12, 6).Value = 96
$ws4.Cells.Item(13, 6).Value = 126
$ws4.Cells.Item(15, 6).Value = 132
$ws4.Cells.Item(17, 6).Value = 255
$ws4.Cells.Item(19, 6).Value = 150
$ws4.Cells.Item(20, 6).Value = 151
$ws4.Cells.Item(21, 6).Value = 1363
$ws4.Cells.Item(22, 6).Value = 517
$ws4.Cells.Item(23, 6).Value = 196
$ws4.Cells.Item(24, 6).Value = 308
$ws4.Cells.Item(26, 6).Value = 1098
$ws4.Cells.Item(27, 6).Value = 2547
$ws4.Cells.Item(29, 6).Value = 1306
$ws4.Cells.Item(34, 6).Value = 219
$ws4.Cells.Item(35, 6).Value = 375
$ws4.Cells.Item(36, 6).Value = 944
$ws4.Cells.Item(39, 6).Value = 779
$ws4.Cells.Item(40, 6).Value = 1091
$ws4.Cells.Item(42, 6).Value = 437
$ws4.Cells.Item(43, 6).Value = 596
$ws4.Cells.Item(44, 6).Value = 762
$ws4.Cells.Item(45, 6).Value = 326
$ws4.Cells.Item(47, 6).Value = 12
$ws4.Cells.Item(48, 6).Value = 216
